# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- Sheet "Rushing" ---
$ws1 = $wb.Worksheets.Item("Rushing")

# Row 3 - T.Huntley
$ws1.Cells.Item(3, 3).Value = 11   # C3
$ws1.Cells.Item(3, 4).Value = 12   # D3
$ws1.Cells.Item(3, 5).Value = 9    # E3

# Row 6 - L.Murray
$ws1.Cells.Item(6, 3).Value = 57   # C6
$ws1.Cells.Item(6, 4).Value = 35   # D6
$ws1.Cells.Item(6, 5).Value = 10   # E6
$ws1.Cells.Item(6, 6).Value = 22   # F6

# Row 7 - D.Freeman
$ws1.Cells.Item(7, 3).Value = 73   # C7
$ws1.Cells.Item(7, 4).Value = 43   # D7
$ws1.Cells.Item(7, 5).Value = 12   # E7
$ws1.Cells.Item(7, 6).Value = 18   # F7

# Row 11 - D.Duvernay
$ws1.Cells.Item(11, 4).Value = 4   # D11

# --- Sheet "Receiving" ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 4 - J.Proche
$ws2.Cells.Item(4, 3).Value = 39   # C4
$ws2.Cells.Item(4, 4).Value = 32   # D4

# Row 6 - M.Brown
$ws2.Cells.Item(6, 3).Value = 96   # C6
$ws2.Cells.Item(6, 4).Value = 73   # D6
$ws2.Cells.Item(6, 5).Value = 42   # E6

# Row 8 - D.Duvernay
$ws2.Cells.Item(8, 3).Value = 44   # C8
$ws2.Cells.Item(8, 4).Value = 32   # D8
$ws2.Cells.Item(8, 7).Value = 8    # G8

# Row 11 - R.Bateman
$ws2.Cells.Item(11, 3).Value = 48  # C11
$ws2.Cells.Item(11, 4).Value = 35  # D11
$ws2.Cells.Item(11, 5).Value = 16  # E11

# Row 13 - M.Andrews
$ws2.Cells.Item(13, 3).Value = 108 # C13
$ws2.Cells.Item(13, 4).Value = 80  # D13
$ws2.Cells.Item(13, 5).Value = 30  # E13
$ws2.Cells.Item(13, 6).Value = 19  # F13

# Row 14 - N.Boyle
$ws2.Cells.Item(14, 3).Value = 2   # C14
$ws2.Cells.Item(14, 4).Value = 1   # D14
$ws2.Cells.Item(14, 7).Value = 1   # G14
$ws2.Cells.Item(14, 8).Value = 1   # H14
